$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GitHubSync")

# Add new D column: orgURL header + URL value
$ws.Range("D1").Value = "orgURL"
$ws.Range("D2").Value = "https://github.com/orgs/dineflesh"

# Rename the username value in row 2
$ws.Range("A2").Value = "bhautik-vasebh1"

# Move the active selection to J10, matching the synced sheet view
$ws.Activate()
[void]$ws.Range("J10").Select()
